$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3:79 down to 4:80
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new data record
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = Get-Date -Year 2022 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100108
$ws.Range("H3").Value = "Tropicales y subtropicales"
$ws.Range("I3").Value = 100108002
$ws.Range("J3").Value = "Mango"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 9500
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9750
$ws.Range("Q3").Value = "$/bandeja 4 kilos"
$ws.Range("R3").Value = "Brasil"
$ws.Range("S3").Value = 2438
$ws.Range("T3").Value = 4
